# Working with Data: Non-Relational Databases - Task complete
# Fill in the "Gross Profit" row (row 6) formulas for FY'14..FY'18 (cols I:L)
# and the "Average Gross Profit" row (row 8) formulas for FY'09..FY'18 (cols C:L),
# replacing the "Put Formula Here" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 "Gross Profit" = Net operating revenues (row 4) - Cost of goods sold (row 5)
# C6:H6 already contain the correct values/formula; extend the pattern to I6:L6.
$ws.Range("I6").Formula = "=I4-I5"
$ws.Range("J6").Formula = "=J4-J5"
$ws.Range("K6").Formula = "=K4-K5"
$ws.Range("L6").Formula = "=L4-L5"

# Row 8 "Average Gross Profit from FY '09 TO FY '18" = AVERAGE(row4, row5) per column
$ws.Range("C8").Formula = "=AVERAGE(C4:C5)"
$ws.Range("D8").Formula = "=AVERAGE(D4:D5)"
$ws.Range("E8").Formula = "=AVERAGE(E4:E5)"
$ws.Range("F8").Formula = "=AVERAGE(F4:F5)"
$ws.Range("G8").Formula = "=AVERAGE(G4:G5)"
$ws.Range("H8").Formula = "=AVERAGE(H4:H5)"
$ws.Range("I8").Formula = "=AVERAGE(I4:I5)"
$ws.Range("J8").Formula = "=AVERAGE(J4:J5)"
$ws.Range("K8").Formula = "=AVERAGE(K4:K5)"
$ws.Range("L8").Formula = "=AVERAGE(L4:L5)"

# Match the cell that was last selected when the workbook was saved.
$ws.Range("D11").Select()
